# "The Last Update 15-03-2024" - refresh NBA team leaders tables.
#
# Numeric-looking "Valor" values are stored as text (shared strings) in
# this workbook, not as real numbers. Assigning a numeric-looking string
# straight to Range.Value makes Excel coerce it into a real number cell,
# which would not match the original text-cell layout. To write a true
# text value (without picking up a new "@" text style), we stage the
# text in a scratch cell via a text formula, copy it, and paste-special
# only the values into the destination - this preserves the destination
# cell's existing (lack of) style while keeping the value as text.

$wb = $excel.ActiveWorkbook

function Set-TextValue($ws, $cellAddr, $text) {
    $scratch = $ws.Range("Z100")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)  # xlPasteValues
    $scratch.ClearContents()
}

# Sheet 1: "Arremessos %"
$ws = $wb.Worksheets.Item("Arremessos %")
$ws.Range("B2").Value = "Indiana Pacers"
Set-TextValue $ws "C2" "50.6"
$ws.Range("B3").Value = "Oklahoma City Thunder"
Set-TextValue $ws "C3" "50.0"
$ws.Range("B4").Value = "Los Angeles Lakers"
Set-TextValue $ws "C4" "49.8"
$ws.Range("B5").Value = "Denver Nuggets"
Set-TextValue $ws "C5" "49.5"
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Phoenix Suns"
Set-TextValue $ws "C6" "49.3"

# Sheet 2: "Diferencial de Pontos"
$ws = $wb.Worksheets.Item("Diferencial de Pontos")
$ws.Range("B2").Value = "Boston Celtics"
Set-TextValue $ws "C2" "+11.2"
$ws.Range("B3").Value = "Oklahoma City Thunder"
Set-TextValue $ws "C3" "+7.6"
$ws.Range("B4").Value = "Minnesota Timberwolves"
Set-TextValue $ws "C4" "+6.4"
$ws.Range("B5").Value = "New Orleans Pelicans"
Set-TextValue $ws "C5" "+4.8"
$ws.Range("B6").Value = "Denver Nuggets"
Set-TextValue $ws "C6" "+4.4"

# Sheet 3: "Pontos"
$ws = $wb.Worksheets.Item("Pontos")
$ws.Range("B2").Value = "Indiana Pacers"
Set-TextValue $ws "C2" "123.1"
$ws.Range("B3").Value = "Boston Celtics"
Set-TextValue $ws "C3" "120.9"
$ws.Range("B4").Value = "Oklahoma City Thunder"
Set-TextValue $ws "C4" "120.8"
$ws.Range("B5").Value = "Milwaukee Bucks"
Set-TextValue $ws "C5" "120.4"
$ws.Range("B6").Value = "Atlanta Hawks"
Set-TextValue $ws "C6" "119.4"

# Sheet 4: "Pontos Permitidos"
$ws = $wb.Worksheets.Item("Pontos Permitidos")
$ws.Range("B2").Value = "Minnesota Timberwolves"
Set-TextValue $ws "C2" "106.7"
$ws.Range("B3").Value = "New York Knicks"
Set-TextValue $ws "C3" "108.1"
$ws.Range("B4").Value = "Orlando Magic"
Set-TextValue $ws "C4" "109.3"
$ws.Range("B5").Value = "Cleveland Cavaliers"
Set-TextValue $ws "C5" "109.4"
$ws.Range("B6").Value = "Boston Celtics"
Set-TextValue $ws "C6" "109.7"

# Sheet 5: "Rebotes"
$ws = $wb.Worksheets.Item("Rebotes")
$ws.Range("B2").Value = "Boston Celtics"
Set-TextValue $ws "C2" "46.9"
$ws.Range("B3").Value = "Golden State Warriors"
Set-TextValue $ws "C3" "46.7"
$ws.Range("B4").Value = "Utah Jazz"
Set-TextValue $ws "C4" "46.1"
$ws.Range("B5").Value = "Houston Rockets"
Set-TextValue $ws "C5" "45.9"
$ws.Range("B6").Value = "New York Knicks"
Set-TextValue $ws "C6" "45.8"

# Sheet 6: "Tocos"
$ws = $wb.Worksheets.Item("Tocos")
$ws.Range("B2").Value = "Oklahoma City Thunder"
Set-TextValue $ws "C2" "6.7"
$ws.Range("B3").Value = "Boston Celtics"
Set-TextValue $ws "C3" "6.6"
$ws.Range("B4").Value = "San Antonio Spurs"
Set-TextValue $ws "C4" "6.3"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Memphis Grizzlies"
Set-TextValue $ws "C5" "6.3"
$ws.Range("B6").Value = "Minnesota Timberwolves"
Set-TextValue $ws "C6" "6.0"
